# ADD results from server
# Update row 2 values on sheets "2025", "2030", "2035" with new data received from server.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" (sheet1.xml) ---
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 44772.20219999998
$ws2025.Range("B2").Value = 66455.8835753044
$ws2025.Range("E2").Value = 161898.8373611681
$ws2025.Range("G2").Value = 42315.16049511277
$ws2025.Range("I2").Value = 368713.9349763304
$ws2025.Range("M2").Value = 117236.9623729033
$ws2025.Range("N2").Value = 44307.00162899461
$ws2025.Range("O2").Value = 69321.42489628839

# --- Sheet "2030" (sheet2.xml) ---
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 0
$ws2030.Range("E2").Value = 168991.0867127558
$ws2030.Range("G2").Value = 0
$ws2030.Range("I2").Value = 186179.998087696
$ws2030.Range("M2").Value = 58342.04889041941
$ws2030.Range("N2").Value = 39676.26863217632
$ws2030.Range("O2").Value = 34021.21142223401

# --- Sheet "2035" (sheet3.xml) ---
$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("E2").Value = 191161.3699534388
$ws2035.Range("G2").Value = 36325.45083361523
$ws2035.Range("I2").Value = 163867.1262511917
$ws2035.Range("L2").Value = 94581.52972418125
$ws2035.Range("M2").Value = 65483.40464893889
$ws2035.Range("N2").Value = 32161.74273064164
$ws2035.Range("O2").Value = 26542.39345168303
